$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'66.110.93"
$ws.Range("E2").Value = "  +2.77%  "

# Row 3
$ws.Range("D3").Value = "'2.961.19"
$ws.Range("E3").Value = "  -0.15%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'574.22"
$ws.Range("E5").Value = "  -0.86%  "

# Row 6
$ws.Range("D6").Value = "'160.58"

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.514"
$ws.Range("E8").Value = "  +0.79%  "

# Row 9
$ws.Range("D9").Value = "'2.957.65"
$ws.Range("E9").Value = "  -0.14%  "

# Row 10
$ws.Range("E10").Value = "  -4.82%  "

# Row 11
$ws.Range("E11").Value = "  -0.86%  "

# Row 12
$ws.Range("D12").Value = "'0.453"
$ws.Range("E12").Value = "  +1.80%  "

# Row 13
$ws.Range("D13").Value = "'0.0000245"
$ws.Range("E13").Value = "  +1.60%  "

# Row 14
$ws.Range("D14").Value = "'34.13"
$ws.Range("E14").Value = "  -0.56%  "

# Row 15
$ws.Range("E15").Value = "  -0.67%  "

# Row 16
$ws.Range("D16").Value = "'66.174.17"
$ws.Range("E16").Value = "  +2.86%  "

# Row 17
$ws.Range("D17").Value = "'3.455.59"
$ws.Range("E17").Value = "  -0.09%  "

# Row 18
$ws.Range("D18").Value = "'6.87"
$ws.Range("E18").Value = "  -0.19%  "

# Row 19
$ws.Range("D19").Value = "'2.961.54"
$ws.Range("E19").Value = "  -0.25%  "

# Row 20
$ws.Range("D20").Value = "'447.28"
$ws.Range("E20").Value = "  +0.77%  "

# Row 21
$ws.Range("D21").Value = "'13.68"
$ws.Range("E21").Value = "  +0.94%  "

# Row 22
$ws.Range("D22").Value = "'0.673"
$ws.Range("E22").Value = "  -0.34%  "

# Row 23
$ws.Range("D23").Value = "'7.13"
$ws.Range("E23").Value = "  -0.99%  "

# Row 24
$ws.Range("D24").Value = "'81.87"
$ws.Range("E24").Value = "  +1.68%  "

# Row 25
$ws.Range("D25").Value = "'2.21"
$ws.Range("E25").Value = "  +0.86%  "

# Row 26
$ws.Range("D26").Value = "'12.13"
$ws.Range("E26").Value = "  -1.26%  "

# Row 27
$ws.Range("E27").Value = "  -0.09%  "

# Row 28
$ws.Range("D28").Value = "'9.96"
$ws.Range("E28").Value = "  -9.30%  "

# Row 29
$ws.Range("D29").Value = "'8.04"
$ws.Range("E29").Value = "  +3.91%  "

# Row 30
$ws.Range("D30").Value = "'2.38"
$ws.Range("E30").Value = "  +8.49%  "

# Row 31
$ws.Range("E31").Value = "  +0.79%  "

# Row 32
$ws.Range("D32").Value = "'0.0₃0984"
$ws.Range("E32").Value = "  -9.88%  "

# Row 33
$ws.Range("D33").Value = "'27.08"
$ws.Range("E33").Value = "  +2.51%  "

# Row 34
$ws.Range("E34").Value = "  -0.81%  "

# Row 35
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.06%  "

# Row 36
$ws.Range("D36").Value = "'0.976"
$ws.Range("E36").Value = "  +0.31%  "

# Row 37
$ws.Range("D37").Value = "'5.68"
$ws.Range("E37").Value = "  +0.95%  "

# Row 38
$ws.Range("D38").Value = "'49.32"
$ws.Range("E38").Value = "  +0.78%  "

# Row 39
$ws.Range("D39").Value = "'2.00"
$ws.Range("E39").Value = "  -4.82%  "

# Row 40
$ws.Range("D40").Value = "'43.32"
$ws.Range("E40").Value = "  -1.44%  "

# Row 41
$ws.Range("D41").Value = "'0.298"
$ws.Range("E41").Value = "  +1.97%  "

# Row 42
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.119"
$ws.Range("E42").Value = "  -0.35%  "

# Row 43
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.81"
$ws.Range("E43").Value = "  -8.26%  "

# Row 44
$ws.Range("D44").Value = "'8.33"
$ws.Range("E44").Value = "  -0.16%  "

# Row 45
$ws.Range("D45").Value = "'381.37"
$ws.Range("E45").Value = "  -1.87%  "

# Row 46
$ws.Range("D46").Value = "'0.0352"
$ws.Range("E46").Value = "  +1.21%  "

# Row 47
$ws.Range("D47").Value = "'2.707.12"
$ws.Range("E47").Value = "  -2.26%  "

# Row 48
$ws.Range("D48").Value = "'131.16"
$ws.Range("E48").Value = "  -2.60%  "

# Row 49
$ws.Range("E49").Value = "  +0.00%  "

# Row 50
$ws.Range("D50").Value = "'0.105"
$ws.Range("E50").Value = "  +0.19%  "

# Row 51
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'2.12"
$ws.Range("E51").Value = "  +4.25%  "
